$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new product-name header cells (I3:W3) ---------------------
$ws.Range("I3").Value = "卡拉蝦原味"
$ws.Range("J3").Value = "卡拉蝦辣味"
$ws.Range("K3").Value = "卡拉魷原味"
$ws.Range("L3").Value = "卡拉魷辣味"
$ws.Range("M3").Value = "卡拉魷芥末"
$ws.Range("N3").Value = "卡拉蟹原味"
$ws.Range("O3").Value = "卡拉蟹辣味"
$ws.Range("P3").Value = "卡拉龍珠原味"
$ws.Range("Q3").Value = "卡拉龍珠辣味"
$ws.Range("R3").Value = "卡拉龍珠芥末"
$ws.Range("S3").Value = "卡拉小卷原味"
$ws.Range("T3").Value = "卡拉小卷芥末"
$ws.Range("U3").Value = "虱目魚薄燒脆片海苔"
$ws.Range("V3").Value = "虱目魚薄燒脆片黑胡椒"
$ws.Range("W3").Value = "虱目魚薄燒脆片蒜香"

# --- Two blank rows above the table, so the used range starts at row 1 --
$ws.Rows.Item(1).RowHeight = $ws.Rows.Item(1).RowHeight
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(2).RowHeight

# --- Explicit ("auto") column widths ------------------------------------
$pad = 5/6
$ws.Columns.Item(1).ColumnWidth = 9 - $pad
$ws.Columns.Item(2).ColumnWidth = 16 - $pad
$ws.Columns.Item(3).ColumnWidth = 15 - $pad
$ws.Columns.Item(4).ColumnWidth = 9 - $pad
$ws.Columns.Item(5).ColumnWidth = 15 - $pad
$ws.Columns.Item(6).ColumnWidth = 16 - $pad
$ws.Columns.Item(7).ColumnWidth = 9 - $pad
$ws.Columns.Item(8).ColumnWidth = 9 - $pad
$ws.Columns.Item(9).ColumnWidth = 10 - $pad
$ws.Columns.Item(10).ColumnWidth = 10 - $pad
$ws.Columns.Item(11).ColumnWidth = 10 - $pad
$ws.Columns.Item(12).ColumnWidth = 10 - $pad
$ws.Columns.Item(13).ColumnWidth = 10 - $pad
$ws.Columns.Item(14).ColumnWidth = 10 - $pad
$ws.Columns.Item(15).ColumnWidth = 10 - $pad
$ws.Columns.Item(16).ColumnWidth = 11 - $pad
$ws.Columns.Item(17).ColumnWidth = 11 - $pad
$ws.Columns.Item(18).ColumnWidth = 11 - $pad
$ws.Columns.Item(19).ColumnWidth = 11 - $pad
$ws.Columns.Item(20).ColumnWidth = 11 - $pad
$ws.Columns.Item(21).ColumnWidth = 14 - $pad
$ws.Columns.Item(22).ColumnWidth = 15 - $pad
$ws.Columns.Item(23).ColumnWidth = 14 - $pad
$ws.Columns.Item(24).ColumnWidth = 7 - $pad
$ws.Columns.Item(25).ColumnWidth = 7 - $pad
$ws.Columns.Item(26).ColumnWidth = 7 - $pad
$ws.Columns.Item(27).ColumnWidth = 7 - $pad
$ws.Columns.Item(28).ColumnWidth = 9 - $pad
$ws.Columns.Item(29).ColumnWidth = 15 - $pad
$ws.Columns.Item(30).ColumnWidth = 15 - $pad
$ws.Columns.Item(31).ColumnWidth = 38 - $pad
$ws.Columns.Item(32).ColumnWidth = 11 - $pad
$ws.Columns.Item(33).ColumnWidth = 7 - $pad
$ws.Columns.Item(34).ColumnWidth = 11 - $pad
$ws.Columns.Item(35).ColumnWidth = 11 - $pad
$ws.Columns.Item(36).ColumnWidth = 7 - $pad
